$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.790.52"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "1.701.56"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'315.01"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "'0.3987"
$ws.Range("E7").Value = "  +2.43%  "
$ws.Range("D8").Value = "'0.4055"
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("B9").Value = "BinanceUSD"
$ws.Range("C9").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D9").Value = "'1.004"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "'53.52"
$ws.Range("E10").Value = "  +1.06%  "
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").Value = "'1.466"
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("D12").Value = "'0.08811"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "'26.25"
$ws.Range("E13").Value = "  +3.44%  "
$ws.Range("D14").Value = "'7.520"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").Value = "'7.965"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "'0.00001341"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "1.827.32"
$ws.Range("E17").Value = "  +8.12%  "
$ws.Range("D18").Value = "'95.65"
$ws.Range("E18").Value = "  -2.90%  "
$ws.Range("D19").Value = "'0.07189"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").Value = "'20.88"
$ws.Range("E20").Value = "  +4.63%  "
$ws.Range("D21").Value = "'7.322"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  +0.91%  "
$ws.Range("D24").Value = "24.778.45"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "'2.376"
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("D26").Value = "'2.879"
$ws.Range("E26").Value = "  -3.43%  "
$ws.Range("D27").Value = "'23.12"
$ws.Range("E27").Value = "  +1.51%  "
$ws.Range("D28").Value = "'6.073"
$ws.Range("E28").Value = "  +16.20%  "
$ws.Range("D29").Value = "'162.04"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "'144.21"
$ws.Range("E30").Value = "  +5.17%  "
$ws.Range("D31").Value = "'8.256"
$ws.Range("E31").Value = "  -5.68%  "
$ws.Range("D32").Value = "2.025.49"
$ws.Range("E32").Value = "  +8.08%  "
$ws.Range("E33").Value = "  +15.44%  "
$ws.Range("D34").Value = "'0.03184"
$ws.Range("E34").Value = "  +8.97%  "
$ws.Range("D35").Value = "'0.08583"
$ws.Range("E35").Value = "  -2.91%  "
$ws.Range("D36").Value = "'7.278"
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("D37").Value = "'1.032"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("E38").Value = "  +3.56%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.8331"
$ws.Range("E39").Value = "  +5.37%  "
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").Value = "'0.09450"
$ws.Range("E40").Value = "  +3.50%  "
$ws.Range("D41").Value = "'10.73"
$ws.Range("E41").Value = "  -0.38%  "
$ws.Range("D42").Value = "'14.19"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("D43").Value = "'1.478"
$ws.Range("E43").Value = "  +1.41%  "
$ws.Range("D44").Value = "'17.50"
$ws.Range("D45").Value = "'2.706"
$ws.Range("E45").Value = "  +4.40%  "
$ws.Range("D46").Value = "'0.7417"
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("D47").Value = "'4.220"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("D48").Value = "'1.373"
$ws.Range("E48").Value = "  +2.59%  "
$ws.Range("D49").Value = "'1.004"
$ws.Range("E49").Value = "  +0.19%  "
$ws.Range("D50").Value = "'0.08382"
$ws.Range("E50").Value = "  +5.18%  "
$ws.Range("D51").Value = "'139.63"
$ws.Range("E51").Value = "  +1.23%  "
